$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests Result ")

# ---------------------------------------------------------------------------
# S001 (rows 15-16): fill in Results / Conclusion text, mark Progression as
# Success.
# ---------------------------------------------------------------------------
$ws.Cells.Item(15, 2).Value = "The error characteristic was done on a scale from 20cm to 600cm. Turned out the sensor can't properly measure anything after 300cm. Between 0cm and 300cm, we can expect an error between 1cm and 5cm. We now need to know if this error is constant in temperature and between different sets or not."
$ws.Cells.Item(15, 2).Font.Italic = $true
$ws.Rows.Item(15).RowHeight = 61.5

$ws.Cells.Item(16, 2).Value = "Further measurement showed that the sensor has a typical error of +-2cm."
$ws.Cells.Item(16, 2).Font.Italic = $true
$ws.Cells.Item(16, 3).Value = "Success"
$ws.Rows.Item(16).RowHeight = 15.4

# ---------------------------------------------------------------------------
# S002 (rows 22-23): no text yet, but Progression moves to "In progress".
# ---------------------------------------------------------------------------
$ws.Cells.Item(22, 2).Font.Italic = $true
$ws.Rows.Item(22).RowHeight = 15.4

$ws.Cells.Item(23, 2).Font.Italic = $true
$ws.Cells.Item(23, 3).Value = "In progress"
$ws.Rows.Item(23).RowHeight = 15.4

# ---------------------------------------------------------------------------
# S003 (rows 29-30): no text yet, but Progression moves to "In progress".
# ---------------------------------------------------------------------------
$ws.Cells.Item(29, 2).Font.Italic = $true
$ws.Rows.Item(29).RowHeight = 15.4

$ws.Cells.Item(30, 2).Font.Italic = $true
$ws.Cells.Item(30, 3).Value = "In progress"
$ws.Rows.Item(30).RowHeight = 15.4

# ---------------------------------------------------------------------------
# S004 (rows 36-37): fill in Results / Conclusion text, mark Progression as
# Success.
# ---------------------------------------------------------------------------
$ws.Cells.Item(36, 2).Value = "The test was done in a temperature-controlled environment with temperature from -15°C to 40°C. Results showed that the sensor is almost not disturbed by the temperature and stay in its +-2cm error from the real distance."
$ws.Cells.Item(36, 2).Font.Italic = $true
$ws.Rows.Item(36).RowHeight = 46.15

$ws.Cells.Item(37, 2).Value = "The sensor pass the test and can be reliably used at various temperatures."
$ws.Cells.Item(37, 2).Font.Italic = $true
$ws.Cells.Item(37, 3).Value = "Success"
$ws.Rows.Item(37).RowHeight = 15.4

# ---------------------------------------------------------------------------
# S005 (rows 43-44): style refresh only, values unchanged (still "Not done").
# ---------------------------------------------------------------------------
$ws.Cells.Item(43, 2).Font.Italic = $true
$ws.Rows.Item(43).RowHeight = 15.4

$ws.Cells.Item(44, 2).Font.Italic = $true
$ws.Rows.Item(44).RowHeight = 15.4

# ---------------------------------------------------------------------------
# S006 (rows 50-51): style refresh only, values unchanged (still "Not done").
# ---------------------------------------------------------------------------
$ws.Cells.Item(50, 2).Font.Italic = $true
$ws.Rows.Item(50).RowHeight = 15.4

$ws.Cells.Item(51, 2).Font.Italic = $true
$ws.Rows.Item(51).RowHeight = 15.4

# ---------------------------------------------------------------------------
# S007 (rows 57-58): style refresh only, values unchanged (still "Not done").
# ---------------------------------------------------------------------------
$ws.Cells.Item(57, 2).Font.Italic = $true
$ws.Rows.Item(57).RowHeight = 15.4

$ws.Cells.Item(58, 2).Font.Italic = $true
$ws.Rows.Item(58).RowHeight = 15.4

# ---------------------------------------------------------------------------
# Sheet view: scroll/selection moved down towards the newly edited rows.
# ---------------------------------------------------------------------------
$ws.Range("C30").Select()

$wb.Save()
